# Auto-generated: update calculated market/profit columns (H-N) across all 8 Job sheets
# per scheduled market-data refresh (Universalis price pull).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1439.1818
$ws.Range("I53").Value = 2150.1428
$ws.Range("J53").Value = 195
$ws.Range("K53").Value = 2150.1428
$ws.Range("L53").Value = 195
$ws.Range("M53").Value = -1513.1428
$ws.Range("N53").Value = -1469
$ws.Range("H62").Value = 2873.6
$ws.Range("I62").Value = 1999.6666
$ws.Range("J62").Value = 4184.5
$ws.Range("K62").Value = 1999.6666
$ws.Range("L62").Value = 4184.5
$ws.Range("M62").Value = -1375.6666
$ws.Range("N62").Value = -5432.5
$ws.Range("H65").Value = 2873.6
$ws.Range("I65").Value = 1999.6666
$ws.Range("J65").Value = 4184.5
$ws.Range("K65").Value = 9998.333000000001
$ws.Range("L65").Value = 20922.5
$ws.Range("M65").Value = -6878.333000000001
$ws.Range("N65").Value = -27162.5
$ws.Range("H116").Value = 15510.5
$ws.Range("I116").Value = 51499
$ws.Range("K116").Value = 51499
$ws.Range("M116").Value = -48057
$ws.Range("H132").Value = 1327.2778
$ws.Range("I132").Value = 1287.7059
$ws.Range("K132").Value = 3863.1177
$ws.Range("M132").Value = -1333.1177
$ws.Range("H137").Value = 2422.1428
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2422.1428
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 7266.428400000001
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -12366.4284
$ws.Range("H138").Value = 2550.238
$ws.Range("I138").Value = 2971.1428
$ws.Range("J138").Value = 2129.3333
$ws.Range("K138").Value = 8913.428400000001
$ws.Range("L138").Value = 6387.999899999999
$ws.Range("M138").Value = -3773.428400000001
$ws.Range("N138").Value = -16667.9999
$ws.Range("H140").Value = 54960.555
$ws.Range("J140").Value = 54960.555
$ws.Range("L140").Value = 54960.555
$ws.Range("N140").Value = -65320.555
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2864.9583
$ws.Range("I61").Value = 2031.1904
$ws.Range("K61").Value = 2031.1904
$ws.Range("M61").Value = -1819.1904
$ws.Range("H74").Value = 1887.7778
$ws.Range("I74").Value = 1623.75
$ws.Range("K74").Value = 1623.75
$ws.Range("M74").Value = -749.75
$ws.Range("H77").Value = 1887.7778
$ws.Range("I77").Value = 1623.75
$ws.Range("K77").Value = 8118.75
$ws.Range("M77").Value = -3750.75
$ws.Range("H136").Value = 2864.9583
$ws.Range("I136").Value = 2031.1904
$ws.Range("K136").Value = 6093.5712
$ws.Range("M136").Value = -3543.5712
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2035.9642
$ws.Range("I105").Value = 2080.261
$ws.Range("J105").Value = 1832.2
$ws.Range("K105").Value = 2080.261
$ws.Range("L105").Value = 1832.2
$ws.Range("M105").Value = -333.261
$ws.Range("N105").Value = -5326.2
$ws.Range("H134").Value = 7705.3105
$ws.Range("I134").Value = 8132.76
$ws.Range("J134").Value = 5033.75
$ws.Range("K134").Value = 24398.28
$ws.Range("L134").Value = 15101.25
$ws.Range("M134").Value = -21863.28
$ws.Range("N134").Value = -20171.25
$ws.Range("H140").Value = 60780
$ws.Range("J140").Value = 60780
$ws.Range("L140").Value = 60780
$ws.Range("N140").Value = -71140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4310.8237
$ws.Range("I31").Value = 1300.5
$ws.Range("J31").Value = 5952.8184
$ws.Range("K31").Value = 1300.5
$ws.Range("L31").Value = 5952.8184
$ws.Range("M31").Value = -1005.5
$ws.Range("N31").Value = -6542.8184
$ws.Range("H34").Value = 4310.8237
$ws.Range("I34").Value = 1300.5
$ws.Range("J34").Value = 5952.8184
$ws.Range("K34").Value = 1300.5
$ws.Range("L34").Value = 5952.8184
$ws.Range("M34").Value = -1098.5
$ws.Range("N34").Value = -6356.8184
$ws.Range("H58").Value = 831
$ws.Range("I58").Value = 802.8333
$ws.Range("K58").Value = 802.8333
$ws.Range("M58").Value = -599.8333
$ws.Range("H94").Value = 1285.4166
$ws.Range("J94").Value = 1258.8572
$ws.Range("L94").Value = 1258.8572
$ws.Range("N94").Value = -2160.8572
$ws.Range("H136").Value = 831
$ws.Range("I136").Value = 802.8333
$ws.Range("K136").Value = 2408.4999
$ws.Range("M136").Value = 141.5001000000002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 299.75
$ws.Range("I8").Value = 299.75
$ws.Range("K8").Value = 899.25
$ws.Range("M8").Value = -760.25
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H109").Value = 1537.3529
$ws.Range("I109").Value = 725.3570999999999
$ws.Range("K109").Value = 2176.0713
$ws.Range("M109").Value = -1136.0713
$ws.Range("H137").Value = 3481.389
$ws.Range("J137").Value = 5836.3335
$ws.Range("L137").Value = 17509.0005
$ws.Range("N137").Value = -27709.0005
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1893.8636
$ws.Range("I102").Value = 2055.3076
$ws.Range("K102").Value = 2055.3076
$ws.Range("M102").Value = -433.3076000000001
$ws.Range("H113").Value = 1504.4445
$ws.Range("I113").Value = 1450
$ws.Range("J113").Value = 1613.3334
$ws.Range("K113").Value = 1450
$ws.Range("L113").Value = 1613.3334
$ws.Range("M113").Value = 720
$ws.Range("N113").Value = -5953.3334
$ws.Range("H128").Value = 25000
$ws.Range("J128").Value = 25000
$ws.Range("L128").Value = 25000
$ws.Range("N128").Value = -34960
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H100").Value = 1667
$ws.Range("I100").Value = 1667.1428
$ws.Range("K100").Value = 1667.1428
$ws.Range("M100").Value = -1126.1428
$ws.Range("H136").Value = 4704.857
$ws.Range("I136").Value = 3428.7144
$ws.Range("J136").Value = 5342.9287
$ws.Range("K136").Value = 10286.1432
$ws.Range("L136").Value = 16028.7861
$ws.Range("M136").Value = -7736.143199999999
$ws.Range("N136").Value = -21128.7861
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 20000
$ws.Range("J32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("N32").Value = -20634
$ws.Range("H100").Value = 386.5
$ws.Range("I100").Value = 298.85715
$ws.Range("K100").Value = 597.7143
$ws.Range("M100").Value = -56.71429999999998
$ws.Range("H136").Value = 2715.2285
$ws.Range("I136").Value = 2260.8
$ws.Range("K136").Value = 6782.400000000001
$ws.Range("M136").Value = -4232.400000000001
